$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 9
$ws1.Range("F6").Value = 160
$ws1.Range("F8").Value = 182
$ws1.Range("F9").Value = 366
$ws1.Range("F10").Value = 478
$ws1.Range("F13").Value = 12141
$ws1.Range("F14").Value = 5451

# Sheet "全部类型" (4th sheet) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 9
$ws4.Range("F8").Value = 160
$ws4.Range("F10").Value = 182
$ws4.Range("F11").Value = 366
$ws4.Range("F12").Value = 478
$ws4.Range("F15").Value = 12141
$ws4.Range("F17").Value = 5451
